$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.228.84"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.59"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.85"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5057"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3929"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09311"
$ws.Range("E9").Value = "  -5.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.143"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.88"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.401"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.89"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.905.55"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.312"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.68"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06615"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.99"
$ws.Range("E20").Value = "  +1.96%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.218"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.263.95"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.321"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.602"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.126.06"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.06"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.25"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.40"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.104"
$ws.Range("E31").Value = "  +2.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1074"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.651"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.612"
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.665"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06657"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02420"
$ws.Range("E37").Value = "  +1.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.249"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.301"
$ws.Range("E39").Value = "  +9.64%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2194"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6449"
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.019"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.50"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6039"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.720"
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.281"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.025"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.21"
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -1.17%  "
